$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule R30's "From" threshold (cell C10) changes from 18 to 1.
$ws.Range("C10").Value = 1
